# harnser-51-sail-measurement.xlsx — "new measurement forms" update
# Updates the luff/foot/leech lengths and the triangle/segment breakdown
# table on the "Form" sheet with the newly-measured dimensions; all the
# dependent formulas (G/I columns, the area totals, and the square-metre
# summary) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form")

# Headline luff / foot / leech measurements (rows 8-10)
$ws.Range("C8").Value = 7423
$ws.Range("C9").Value = 2598
$ws.Range("C10").Value = 7836

# Positive-area triangulation table (rows 13-20)
# Row 13
$ws.Range("C13").Value = 7836
$ws.Range("E13").Value = 2459

# Row 14
$ws.Range("C14").Value = 7423
$ws.Range("E14").Value = 112

# Row 15
$ws.Range("B15").Value = "Triangle"
$ws.Range("C15").Value = 2598
$ws.Range("E15").Value = 124

# Row 16
$ws.Range("C16").Value = 7836
$ws.Range("E16").Value = 1019

# Row 17
$ws.Range("B17").Value = "Triangle"
$ws.Range("C17").Value = 1625
$ws.Range("E17").Value = 645

# Row 18
$ws.Range("B18").Value = "Segment"
$ws.Range("C18").Value = 6641
$ws.Range("E18").Value = 55

# Row 19 - now blank (no third segment recorded for this panel)
$ws.Range("B19").Value = "Triangle"
$ws.Range("C19").Value = ""
$ws.Range("E19").Value = ""

# Row 20 - now blank
$ws.Range("C20").Value = ""
$ws.Range("E20").Value = ""

# Move the active selection/scroll position to reflect where the
# measurer was last working on the form
$null = $ws.Range("N10").Select()
